# "forgot to add the links!" -- add the three missing Amazon product
# links into column F (the "Link" column) for the rows that were
# missing them (elastic, buckle clips, slide buckles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value  = "https://www.amazon.co.uk/dp/9553548296/?coliid=I1WUO41FZR9AN4&colid=4EYVC5LB6BFN&psc=0&ref_=lv_ov_lig_dp_it"
$ws.Range("F9").Value  = "https://www.amazon.co.uk/dp/B00MNKUIT8/?coliid=IG139WPE5GBEE&colid=4EYVC5LB6BFN&psc=1&ref_=lv_ov_lig_dp_it"
$ws.Range("F10").Value = "https://www.amazon.co.uk/dp/B00J3YAK2M/?coliid=IZ7CDRDN889D&colid=4EYVC5LB6BFN&psc=0&ref_=lv_ov_lig_dp_it"

# Reflect where the user ended up on-screen after pasting the links in:
# scrolled right so column E is the leftmost visible column, with F12
# selected as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F12").Select()
